$d = $word.ActiveDocument

# 1. "Productive Systems I" -> "Productive Systems" (Heading3)
$d.Content.Find.Execute("Productive Systems I", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Productive Systems", 2)

# 2. Ativação date update
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2025", 2)

# 3. Portuguese "Programa" detailed paragraph
$oldPt = "1 – Introdução aos Sistemas ProdutivosProdução na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da ProduçãoPapel da função produção. Objetivos de desempenho. 3 – Estratégia de ProduçãoIntrodução. Prioridade de objetivos de desempenho. Áreas de decisão da estratégia de operações.4 – Projeto em Gestão de ProduçãoDefinição de projeto. Principais aspectos de um projeto. Tipos de processos em manufatura e serviços. 5 – Projeto de Produtos e ServiçosIntrodução. Geração de conceito. Triagem de conceito. Avaliação e melhoria do projeto. Protótipo e projeto final.6 – Projeto da Rede de Operações ProdutivasPerspectiva da rede. Integração Vertical. Localização da capacidade. Gestão da capacidade produtiva de longo prazo.7 – Arranjo Físico e FluxoProcedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico."
$newPt = "1 – Introdução aos Sistemas Produtivos; Produção na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da Produção; Papel da função produção. Objetivos de desempenho.3 – Tipos de Manufatura; Tipos básicos de Manufatura.4 – Arranjo Físico e Fluxo; Procedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.5 – Organização do Trabalho e Métodos;Técnicas de organização e métodos de trabalho6 - Introdução ao Planejamento e Controle de Produção.Conceituação do PCP; conciliação de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.7 - Introdução à qualidade e a tecnologia de processo Importância; visões; princípios de administração da qualidade total."
$d.Content.Find.Execute($oldPt, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newPt, 2)

# 4. English (italic) "Programa" detailed paragraph
$oldEn = "1 - Introduction to Productive SystemsProduction in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 - Strategic Role of ProductionRole of production function. Performance Objectives.3 - Production StrategyIntroduction. Priority of performance goals. Operations strategy decision areas.4 - Project in Production ManagementDefinition of project. Main aspects of a project. Types of processes in manufacturing and services.5 - Product and Service ProjectIntroduction. Concept generation. Concept screening. Evaluation and improvement of the project. Prototype and final design.6 - Production Operations Network ProjectNetwork perspective. Vertical integration. Location of capacity. Management of long-term productive capacity.7 – Layout and FlowLayout and Physical Arrangement Procedure. Basic types of physical arrangement. Design of layout and physical arrangement."
$newEn = "1 – Introduction to Production Systems;Production in the Organization. Inputs, Transformation Processes and Outputs. Types of Production Operations. Production management activities.2 – Strategic Role of Production;Role of the production function. Performance objectives.3 – Types of Manufacturing;Basic types of Manufacturing.4 – Physical Arrangement and Flow;Physical Arrangement Procedure. Basic types of physical arrangement. Physical arrangement design.5 – Work Organization and Methods;Organization techniques and work methods6 - Introduction to Production Planning and Control.Conceptualization of the PCP; reconciliation of supply and demand; nature of supply and demand; PCP activities; volume-variety effect in PCP.7 - Introduction to quality and process technologyImportance; visions; total quality management principles."
$d.Content.Find.Execute($oldEn, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newEn, 2)
